{"js": "// Update the ProfessorMenu bullet describing the course-creation form\n// fields: drop \"\u0646\u0627\u0645 \u0627\u0633\u062a\u0627\u062f\" from the list and append \"\u0648 \u0622\u06cc\u062f\u06cc \u0627\u0633\u062a\u0627\u062f\" at\n// the end of the parenthetical, matching the StudentMenu/ProfessorMenu\n// wording update described in the commit.\n\nconst oldText =\n  \"\u0627\u0633\u062a\u0627\u062f \u0627\u0637\u0644\u0627\u0639\u0627\u062a \u062f\u0631\u0633 \u0631\u0627 (\u0639\u0646\u0648\u0627\u0646 \u062f\u0631\u0633\u060c \u0646\u0627\u0645 \u0627\u0633\u062a\u0627\u062f\u060c \u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633 \u0648 \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647) \u0648\u0627\u0631\u062f \u06a9\u0631\u062f\u0647 \u0648 \u06af\u0632\u06cc\u0646\u0647 \u062a\u06a9\u0645\u06cc\u0644 \u0627\u06cc\u062c\u0627\u062f \u062f\u0631\u0633 \u0631\u0627 \u0627\u0646\u062a\u062e\u0627\u0628 \u0645\u06cc \u06a9\u0646\u062f\";\nconst newText =\n  \"\u0627\u0633\u062a\u0627\u062f \u0627\u0637\u0644\u0627\u0639\u0627\u062a \u062f\u0631\u0633 \u0631\u0627 (\u0639\u0646\u0648\u0627\u0646 \u062f\u0631\u0633\u060c \u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633\u060c \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647 \u0648 \u0622\u06cc\u062f\u06cc \u0627\u0633\u062a\u0627\u062f) \u0648\u0627\u0631\u062f \u06a9\u0631\u062f\u0647 \u0648 \u06af\u0632\u06cc\u0646\u0647 \u062a\u06a9\u0645\u06cc\u0644 \u0627\u06cc\u062c\u0627\u062f \u062f\u0631\u0633 \u0631\u0627 \u0627\u0646\u062a\u062e\u0627\u0628 \u0645\u06cc \u06a9\u0646\u062f\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(newText, \"Replace\");\n} else {\n  // Fallback: locate just the changed middle fragment in case the\n  // surrounding text was already normalized differently, and rebuild\n  // the full sentence around it.\n  const middleOld = \"\u0646\u0627\u0645 \u0627\u0633\u062a\u0627\u062f\u060c \u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633 \u0648 \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647\";\n  const middleNew = \"\u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633\u060c \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647 \u0648 \u0622\u06cc\u062f\u06cc \u0627\u0633\u062a\u0627\u062f\";\n  const middleResults = context.document.body.search(middleOld, {\n    matchCase: true,\n  });\n  await context.sync();\n  if (middleResults.items.length > 0) {\n    middleResults.items[0].insertText(middleNew, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the ProfessorMenu bullet describing the course-creation form\n# fields: drop \"\u0646\u0627\u0645 \u0627\u0633\u062a\u0627\u062f\" from the parenthetical list and append\n# \"\u0648 \u0622\u06cc\u062f\u06cc \u0627\u0633\u062a\u0627\u062f\" at the end, matching the StudentMenu/ProfessorMenu\n# wording update described in the commit.\n\n$d = $word.ActiveDocument\n\n$oldText = \"\u0627\u0633\u062a\u0627\u062f \u0627\u0637\u0644\u0627\u0639\u0627\u062a \u062f\u0631\u0633 \u0631\u0627 (\u0639\u0646\u0648\u0627\u0646 \u062f\u0631\u0633\u060c \u0646\u0627\u0645 \u0627\u0633\u062a\u0627\u062f\u060c \u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633 \u0648 \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647) \u0648\u0627\u0631\u062f \u06a9\u0631\u062f\u0647 \u0648 \u06af\u0632\u06cc\u0646\u0647 \u062a\u06a9\u0645\u06cc\u0644 \u0627\u06cc\u062c\u0627\u062f \u062f\u0631\u0633 \u0631\u0627 \u0627\u0646\u062a\u062e\u0627\u0628 \u0645\u06cc \u06a9\u0646\u062f\"\n$newText = \"\u0627\u0633\u062a\u0627\u062f \u0627\u0637\u0644\u0627\u0639\u0627\u062a \u062f\u0631\u0633 \u0631\u0627 (\u0639\u0646\u0648\u0627\u0646 \u062f\u0631\u0633\u060c \u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633\u060c \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647 \u0648 \u0622\u06cc\u062f\u06cc \u0627\u0633\u062a\u0627\u062f) \u0648\u0627\u0631\u062f \u06a9\u0631\u062f\u0647 \u0648 \u06af\u0632\u06cc\u0646\u0647 \u062a\u06a9\u0645\u06cc\u0644 \u0627\u06cc\u062c\u0627\u062f \u062f\u0631\u0633 \u0631\u0627 \u0627\u0646\u062a\u062e\u0627\u0628 \u0645\u06cc \u06a9\u0646\u062f\"\n\n$range = $d.Content\n$found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\nWrite-Host \"Replaced main sentence: $found\"\n\nif (-not $found) {\n    # Fallback: only the changed middle fragment, in case the exact\n    # full-sentence text was not found verbatim.\n    $middleOld = \"\u0646\u0627\u0645 \u0627\u0633\u062a\u0627\u062f\u060c \u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633 \u0648 \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647\"\n    $middleNew = \"\u0633\u0631\u062a\u0631\u0645 \u062f\u0631\u0633\u060c \u0646\u0627\u0645 \u062f\u0627\u0646\u0634\u06a9\u062f\u0647 \u0648 \u0622\u06cc\u062f\u06cc \u0627\u0633\u062a\u0627\u062f\"\n    $range2 = $d.Content\n    $found2 = $range2.Find.Execute($middleOld, $true, $false, $false, $false, $false, $true, 1, $false, $middleNew, 2)\n    Write-Host \"Replaced fallback fragment: $found2\"\n}\n"}
